$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 16, shifting rows 16-93 down to 17-94.
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with the new data record.
$ws.Cells.Item(16, 1).Value = 9
$ws.Cells.Item(16, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(16, 3).Value = "Metropolitana"
$ws.Cells.Item(16, 4).Value = 44547
$ws.Cells.Item(16, 5).Value = 13
$ws.Cells.Item(16, 6).Value = 100112022
$ws.Cells.Item(16, 7).Value = "Arveja Verde"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 43
$ws.Cells.Item(16, 11).Value = 11000
$ws.Cells.Item(16, 12).Value = 12000
$ws.Cells.Item(16, 13).Value = 11512
$ws.Cells.Item(16, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(16, 15).Value = "Carahue"
$ws.Cells.Item(16, 16).Value = 460
$ws.Cells.Item(16, 17).Value = 25
$ws.Cells.Item(16, 18).Value = "Hortaliza"
